$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.272.68'
$ws.Range("E2").Value = '  -1.94%  '

$ws.Range("D3").Value = '3.157.94'
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("E4").Value = '  -0.09%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '590.99'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.29%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '138.92'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -4.08%  '

$ws.Range("E7").Value = '  -0.60%  '

$ws.Range("D8").Value = '3.155.57'
$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("E9").Value = '  -1.34%  '

$ws.Range("E10").Value = '  -2.55%  '

$ws.Range("E11").Value = '  -2.02%  '

$ws.Range("E12").Value = '  -2.65%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000246'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.58%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '34.26'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.85%  '

$ws.Range("D15").Value = '3.678.14'

$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("D17").Value = '3.153.75'
$ws.Range("E17").Value = '  -0.48%  '

$ws.Range("D18").Value = '63.221.32'

$ws.Range("E19").Value = '  -3.02%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '476.29'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.45%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '14.11'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -4.05%  '

$ws.Range("E22").Value = '  -1.91%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '7.73'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '84.62'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -3.62%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '13.02'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -3.78%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("E27").Value = '  -1.98%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.11'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.07%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '8.02'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -5.09%  '

$ws.Range("E30").Value = '  +1.51%  '

$ws.Range("E31").Value = '  -0.02%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '26.95'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.65%  '

$ws.Range("E33").Value = '  -5.35%  '

$ws.Range("E34").Value = '  -5.56%  '

$ws.Range("E35").Value = '  -3.38%  '

$ws.Range("E36").Value = '  -4.10%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '52.54'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.87%  '

$ws.Range("D38").Value = '0.0₃0702'
$ws.Range("E38").Value = '  -8.03%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0390'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.90%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '422.65'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -4.70%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -9.42%  '

$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("D43").Value = '2.938.51'
$ws.Range("E43").Value = '  +2.25%  '

$ws.Range("E44").Value = '  -6.15%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.264'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.21%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.14'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.87%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '25.56'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.33%  '

$ws.Range("E49").Value = '  -0.67%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -9.65%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '121.12'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.60%  '

